# Update data values on the "north" and "south" sheets.
# The diff only changes the numeric values in columns C:F for rows 2-5
# on both sheets (no structural/formatting changes).

$wb = $excel.ActiveWorkbook

$north = $wb.Worksheets.Item("north")
$south = $wb.Worksheets.Item("south")

# north sheet (sheet1) new values for rows 2-5, columns C:F
$north.Cells.Item(2, 3).Value = 23
$north.Cells.Item(2, 4).Value = 8
$north.Cells.Item(2, 5).Value = 6
$north.Cells.Item(2, 6).Value = 31

$north.Cells.Item(3, 3).Value = 7
$north.Cells.Item(3, 4).Value = 21
$north.Cells.Item(3, 5).Value = 9
$north.Cells.Item(3, 6).Value = 22

$north.Cells.Item(4, 3).Value = 13
$north.Cells.Item(4, 4).Value = 19
$north.Cells.Item(4, 5).Value = 4
$north.Cells.Item(4, 6).Value = 15

$north.Cells.Item(5, 3).Value = 32
$north.Cells.Item(5, 4).Value = 20
$north.Cells.Item(5, 5).Value = 2
$north.Cells.Item(5, 6).Value = 11

# south sheet (sheet2) new values for rows 2-5, columns C:F
$south.Cells.Item(2, 3).Value = 14
$south.Cells.Item(2, 4).Value = 1
$south.Cells.Item(2, 5).Value = 12
$south.Cells.Item(2, 6).Value = 25

$south.Cells.Item(3, 3).Value = 29
$south.Cells.Item(3, 4).Value = 18
$south.Cells.Item(3, 5).Value = 3
$south.Cells.Item(3, 6).Value = 24

$south.Cells.Item(4, 3).Value = 10
$south.Cells.Item(4, 4).Value = 5
$south.Cells.Item(4, 5).Value = 16
$south.Cells.Item(4, 6).Value = 17

$south.Cells.Item(5, 3).Value = 28
$south.Cells.Item(5, 4).Value = 27
$south.Cells.Item(5, 5).Value = 26
$south.Cells.Item(5, 6).Value = 30
